# Apply the StoryCards.xlsx edit: update status of row 21 & 22, add effort/date
# info to row 21, and move the active selection to B22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 21 (Story ID 14): status "in Arbeit" -> "fertig"
$ws.Range("B21").Value = "fertig"

# Row 21: add geschätzter Aufwand (K), tatsächlicher Aufwand (L), Fertigstelldatum (M)
$ws.Range("K21").Value = "2h"
$ws.Range("L21").Value = "3h"
# Give M21 the same date format as the neighbouring date column (E21) by
# copying its formatting (xlPasteFormats = -4122), then set the date value
# (40823 = 2011-10-07 as an OLE Automation date serial).
$ws.Range("E21").Copy()
$ws.Range("M21").PasteSpecial(-4122)
$ws.Range("M21").Value = 40823

# Row 22 (Story ID 15): status "jungfräulich" -> "in Arbeit"
$ws.Range("B22").Value = "in Arbeit"

# Update the active selection shown in the sheet view
$ws.Range("B22").Select()
